$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 106 (shifts existing rows 106-211 down to 107-212)
$ws.Rows.Item(106).Insert()

# Populate the new row 106 with the new data point
$ws.Range("A106").Value = 8
$ws.Range("B106").Value = "Terminal La Palmera de La Serena"
$ws.Range("C106").Value = "Coquimbo"
$ws.Range("D106").Value = 45271
$ws.Range("E106").Value = 4
$ws.Range("F106").Value = 100114007
$ws.Range("G106").Value = "Jengibre"
$ws.Range("H106").Value = "Sin especificar"
$ws.Range("I106").Value = "Primera"
$ws.Range("J106").Value = 340
$ws.Range("K106").Value = 22000
$ws.Range("L106").Value = 23000
$ws.Range("M106").Value = 22500
$ws.Range("N106").Value = "$/caja 13 kilos"
$ws.Range("O106").Value = "Perú"
$ws.Range("P106").Value = 1731
$ws.Range("Q106").Value = 13
$ws.Range("R106").Value = "Hortaliza"
